# Finnish_data.xlsx maintenance edit:
#  - Move the "comms" (col M) free-text notes into "comms_internal" (col S),
#    which is the column actually used for internal remarks.
#  - For rows whose valency pattern was plain "TR" (X/Y both bare
#    transitive), split it into the real case frame NOM (subject) +
#    GEN~NOM (object alternation) in columns I/J.
#  - For the two "*" (unspecified) rows, mirror that marker into column J too.
#  - Move the active-cell selection from M9 to M5 (it now points at an
#    empty "comms" cell instead of the moved note).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Finnish")

# --- rows where I/J used to be "TR"/<blank> -> becomes "NOM"/"GEN~NOM" ---
$ijRows = @(5,9,10,12,14,17,20,21,22,27,28,30,34,37,41,43,45,50,51,56,64,69,70,71,72,73,76,77,86,87,97,101,103,104,106,109,122,125)
foreach ($r in $ijRows) {
    $ws.Cells.Item($r, 9).Value  = "NOM"      # column I
    $ws.Cells.Item($r, 10).Value = "GEN~NOM"  # column J
}

# --- rows where I = "*" and J was blank -> J becomes "*" as well ---
$starRows = @(25,44)
foreach ($r in $starRows) {
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 9).Value2
}

# --- move every non-empty "comms" (col M) value over to "comms_internal" (col S) ---
$mRows = @(2,3,6,7,13,15,19,20,22,23,25,32,35,36,39,43,44,45,47,48,51,52,54,57,59,60,62,65,66,67,68,73,81,82,84,85,86,87,89,91,94,98,99,110,112,113,114,116,118,120,121,123,125,126,130,131)

# Normalise every comms_internal cell (col S, rows 2-131) to the same format
# used by the comms column before clearing it out, then fill in the moved text.
[void]$ws.Range("M2").Copy()
[void]$ws.Range("S2:S131").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($r in $mRows) {
    $note = $ws.Cells.Item($r, 13).Value2
    $ws.Cells.Item($r, 19).Value = $note
    $ws.Cells.Item($r, 13).Value = ""
}

# --- move the active selection from M9 to M5 ---
[void]$ws.Range("M5").Select()
